$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel
# auto-converting a numeric-looking string ("1.003", "0.5867", ...)
# into a real number. We briefly mark the cell as Text, assign the
# value, then restore the cell style so no extra formatting sticks.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '27.502.51'
$ws.Range("E2").Value = '  -2.15%  '
$ws.Range("D3").Value = '1.748.49'
$ws.Range("E3").Value = '  -3.03%  '
Set-TextValue "D4" '1.003'
$ws.Range("E4").Value = '  -0.13%  '
Set-TextValue "D5" '324.13'
$ws.Range("E5").Value = '  -0.09%  '
Set-TextValue "D6" '1.002'
$ws.Range("E6").Value = '  -0.08%  '
Set-TextValue "D7" '0.4415'
$ws.Range("E7").Value = '  +3.17%  '
Set-TextValue "D8" '0.3608'
$ws.Range("E8").Value = '  -0.56%  '
Set-TextValue "D9" '0.07468'
$ws.Range("E9").Value = '  -1.38%  '
Set-TextValue "D10" '42.25'
$ws.Range("E10").Value = '  -5.78%  '
$ws.Range("E11").Value = '  -2.67%  '
Set-TextValue "D12" '1.003'
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("E13").Value = '  -4.76%  '
Set-TextValue "D14" '6.024'
$ws.Range("E14").Value = '  -3.35%  '
Set-TextValue "D15" '7.132'
$ws.Range("E15").Value = '  -3.50%  '
$ws.Range("D16").Value = '1.751.10'
$ws.Range("E16").Value = '  -3.99%  '
Set-TextValue "D17" '92.13'
$ws.Range("E17").Value = '  -1.15%  '
Set-TextValue "D18" '0.00001057'
$ws.Range("E18").Value = '  -1.12%  '
Set-TextValue "D19" '0.06408'
$ws.Range("E19").Value = '  +0.60%  '
$ws.Range("E20").Value = '  +0.06%  '
Set-TextValue "D21" '16.83'
$ws.Range("E21").Value = '  -2.49%  '
$ws.Range("E22").Value = '  -4.59%  '
$ws.Range("D23").Value = '27.557.50'
$ws.Range("E23").Value = '  -2.07%  '
Set-TextValue "D24" '11.17'
$ws.Range("E24").Value = '  -2.54%  '
Set-TextValue "D25" '2.097'
$ws.Range("E25").Value = '  -2.51%  '
Set-TextValue "D26" '161.36'
$ws.Range("E26").Value = '  +0.66%  '
Set-TextValue "D27" '20.36'
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("D28").Value = '1.949.61'
$ws.Range("E28").Value = '  -3.87%  '
Set-TextValue "D29" '2.102'
$ws.Range("E29").Value = '  -5.83%  '
Set-TextValue "D30" '124.48'
$ws.Range("E30").Value = '  -3.74%  '
$ws.Range("E31").Value = '  -8.84%  '
Set-TextValue "D32" '3.652'
$ws.Range("E32").Value = '  +3.18%  '
Set-TextValue "D33" '0.08946'
Set-TextValue "D34" '5.516'
$ws.Range("E34").Value = '  -6.62%  '
$ws.Range("E35").Value = '  -6.83%  '
Set-TextValue "D36" '0.02311'
$ws.Range("E36").Value = '  -2.90%  '
Set-TextValue "D37" '0.2083'
$ws.Range("E37").Value = '  -2.15%  '
Set-TextValue "D38" '0.6329'
$ws.Range("E38").Value = '  -2.69%  '
Set-TextValue "D39" '0.05966'
$ws.Range("E39").Value = '  -2.68%  '
Set-TextValue "D40" '4.926'
$ws.Range("E40").Value = '  -4.13%  '
Set-TextValue "D41" '1.202'
$ws.Range("E41").Value = '  +0.43%  '
$ws.Range("E42").Value = '  -0.09%  '
Set-TextValue "D43" '1.383'
$ws.Range("E43").Value = '  -3.04%  '
Set-TextValue "D44" '7.746'
$ws.Range("E44").Value = '  -2.91%  '
Set-TextValue "D45" '13.24'
$ws.Range("E45").Value = '  -3.11%  '
Set-TextValue "D46" '3.710'
$ws.Range("E46").Value = '  -0.36%  '
Set-TextValue "D47" '0.5867'
$ws.Range("E47").Value = '  -2.49%  '
Set-TextValue "D48" '121.11'
$ws.Range("E48").Value = '  -3.23%  '
Set-TextValue "D49" '1.942'
$ws.Range("E49").Value = '  -2.55%  '
Set-TextValue "D50" '1.149'
$ws.Range("E50").Value = '  -1.27%  '
Set-TextValue "D51" '0.06862'
$ws.Range("E51").Value = '  -1.66%  '
